$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 78, shifting the existing rows (and the
# "800"/"900" category blocks below) down by one.
$ws.Range("A78").EntireRow.Insert()

# Populate the newly inserted row with the new log-message entry.
$ws.Range("A78").Value = "The custom property name %s contains whitespace, replacing it with underscores"
$ws.Range("B78").Value = "Possible issue in original STIX 1.x content"
$ws.Range("C78").Value = 624
$ws.Range("D78").Value = "warn"
$ws.Range("E78").Value = "convert_custom_properties"

# Match the view state left behind in the saved workbook.
$ws.Application.ActiveWindow.ScrollRow = 58
$ws.Range("F78").Select()
